# Update the 江西-漫展信息 workbook: remove the 2024-03-02 "meeting" event row
# (all later rows shift up by one), renumber the index column (A), and bump
# the "想去人数" (F) counters that were refreshed when the page was re-scraped.
#
# This same transformation is applied identically to the two data sheets:
#   展览  (exhibitions)
#   全部类型 (all types)
# (演出 and 本地生活 only contain the header row and are untouched.)

$wb = $excel.ActiveWorkbook

# New F-column ("想去人数") values, keyed by the *final* row number (1-based,
# row 1 is the header) after the obsolete row has been removed and every
# later row has shifted up by one.
$fUpdates = @{
    2  = 59
    4  = 37
    6  = 2791
    8  = 1800
    9  = 188
    11 = 639
    14 = 169
    15 = 85
    16 = 88
    17 = 24
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Remove the old row 2 (2024-03-02 南昌·meeting动漫游戏嘉年华).
    # Everything below shifts up by one row automatically.
    $ws.Rows("2:2").Delete()

    # Column A holds a plain sequential index (row number - 1); after the
    # deletion it still contains the stale values copied up from below, so
    # renumber it for the surviving data rows (now rows 2-17).
    for ($r = 2; $r -le 17; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }

    # Apply the refreshed "want to go" counts.
    foreach ($r in $fUpdates.Keys) {
        $ws.Cells.Item($r, 6).Value = $fUpdates[$r]
    }
}
